# Auto-generated edit script applying scheduled-runner market-data refresh
# to the Sagittarius_Profits workbook (8 sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each touched row holds static (non-formula) cached market values in columns H:N:
#   H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#   K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 29.5
$ws.Range("I33").Value = 29.5
$ws.Range("K33").Value = 29.5
$ws.Range("M33").Value = 199.5
$ws.Range("H80").Value = 845.1429000000001
$ws.Range("J80").Value = 1047.8
$ws.Range("L80").Value = 3143.4
$ws.Range("N80").Value = -5139.4
$ws.Range("H83").Value = 845.1429000000001
$ws.Range("J83").Value = 1047.8
$ws.Range("L83").Value = 9430.199999999999
$ws.Range("N83").Value = -19414.2
$ws.Range("H132").Value = 4004.5
$ws.Range("I132").Value = 4004.5
$ws.Range("K132").Value = 12013.5
$ws.Range("M132").Value = -9483.5
$ws.Range("H141").Value = 2748.75
$ws.Range("I141").Value = 2748.75
$ws.Range("K141").Value = 8246.25
$ws.Range("M141").Value = -3066.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2117.7368
$ws.Range("I2").Value = 1648.7
$ws.Range("J2").Value = 2638.889
$ws.Range("K2").Value = 1648.7
$ws.Range("L2").Value = 2638.889
$ws.Range("M2").Value = -1535.7
$ws.Range("N2").Value = -2864.889
$ws.Range("H32").Value = 3351075.8
$ws.Range("I32").Value = 3503637.2
$ws.Range("K32").Value = 3503637.2
$ws.Range("M32").Value = -3503350.2
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H74").Value = 1767.7778
$ws.Range("I74").Value = 1485
$ws.Range("J74").Value = 2333.3333
$ws.Range("K74").Value = 1485
$ws.Range("L74").Value = 2333.3333
$ws.Range("M74").Value = -611
$ws.Range("N74").Value = -4081.3333
$ws.Range("H77").Value = 1767.7778
$ws.Range("I77").Value = 1485
$ws.Range("J77").Value = 2333.3333
$ws.Range("K77").Value = 7425
$ws.Range("L77").Value = 11666.6665
$ws.Range("M77").Value = -3057
$ws.Range("N77").Value = -20402.6665
$ws.Range("H116").Value = 2117.7368
$ws.Range("I116").Value = 1648.7
$ws.Range("J116").Value = 2638.889
$ws.Range("K116").Value = 1648.7
$ws.Range("L116").Value = 2638.889
$ws.Range("M116").Value = 645.3
$ws.Range("N116").Value = -7226.889
$ws.Range("H122").Value = 4110.3335
$ws.Range("I122").Value = 3332
$ws.Range("K122").Value = 9996
$ws.Range("M122").Value = -7546
$ws.Range("H132").Value = 1499
$ws.Range("J132").Value = 1499
$ws.Range("L132").Value = 4497
$ws.Range("N132").Value = -9557

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2117.7368
$ws.Range("I3").Value = 1648.7
$ws.Range("J3").Value = 2638.889
$ws.Range("K3").Value = 1648.7
$ws.Range("L3").Value = 2638.889
$ws.Range("M3").Value = -1534.7
$ws.Range("N3").Value = -2866.889
$ws.Range("H99").Value = 1987.6666
$ws.Range("I99").Value = 1995
$ws.Range("J99").Value = 1980.3334
$ws.Range("K99").Value = 1995
$ws.Range("L99").Value = 1980.3334
$ws.Range("M99").Value = -497
$ws.Range("N99").Value = -4976.3334
$ws.Range("H105").Value = 2082
$ws.Range("I105").Value = 2225
$ws.Range("K105").Value = 2225
$ws.Range("M105").Value = -478

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1525.9231
$ws.Range("I7").Value = 492.85715
$ws.Range("J7").Value = 2731.1667
$ws.Range("K7").Value = 492.85715
$ws.Range("L7").Value = 2731.1667
$ws.Range("M7").Value = -379.85715
$ws.Range("N7").Value = -2957.1667
$ws.Range("H31").Value = 1687.8
$ws.Range("I31").Value = 1813
$ws.Range("J31").Value = 1500
$ws.Range("K31").Value = 1813
$ws.Range("L31").Value = 1500
$ws.Range("M31").Value = -1518
$ws.Range("N31").Value = -2090
$ws.Range("H34").Value = 1687.8
$ws.Range("I34").Value = 1813
$ws.Range("J34").Value = 1500
$ws.Range("K34").Value = 1813
$ws.Range("L34").Value = 1500
$ws.Range("M34").Value = -1611
$ws.Range("N34").Value = -1904
$ws.Range("H86").Value = 5510.4
$ws.Range("I86").Value = 4027.4285
$ws.Range("J86").Value = 8970.666999999999
$ws.Range("K86").Value = 4027.4285
$ws.Range("L86").Value = 8970.666999999999
$ws.Range("M86").Value = -2904.4285
$ws.Range("N86").Value = -11216.667
$ws.Range("H89").Value = 5510.4
$ws.Range("I89").Value = 4027.4285
$ws.Range("J89").Value = 8970.666999999999
$ws.Range("K89").Value = 20137.1425
$ws.Range("L89").Value = 44853.335
$ws.Range("M89").Value = -14521.1425
$ws.Range("N89").Value = -56085.335
$ws.Range("H99").Value = 1773.2
$ws.Range("J99").Value = 1891.5
$ws.Range("L99").Value = 1891.5
$ws.Range("N99").Value = -4887.5
$ws.Range("H105").Value = 2896.2727
$ws.Range("I105").Value = 2195.2
$ws.Range("J105").Value = 3480.5
$ws.Range("K105").Value = 2195.2
$ws.Range("L105").Value = 3480.5
$ws.Range("M105").Value = -448.1999999999998
$ws.Range("N105").Value = -6974.5
$ws.Range("H122").Value = 2585
$ws.Range("I122").Value = 3311.6667
$ws.Range("J122").Value = 1495
$ws.Range("K122").Value = 9935.000100000001
$ws.Range("L122").Value = 4485
$ws.Range("M122").Value = -7485.000100000001
$ws.Range("N122").Value = -9385
$ws.Range("H126").Value = 1773.2
$ws.Range("J126").Value = 1891.5
$ws.Range("L126").Value = 5674.5
$ws.Range("N126").Value = -10614.5
$ws.Range("H134").Value = 2522.0625
$ws.Range("I134").Value = 2490.2
$ws.Range("K134").Value = 7470.599999999999
$ws.Range("M134").Value = -4935.599999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1499.5
$ws.Range("I75").Value = 1499
$ws.Range("J75").Value = 1500
$ws.Range("K75").Value = 4497
$ws.Range("L75").Value = 4500
$ws.Range("M75").Value = -3499
$ws.Range("N75").Value = -6496
$ws.Range("H78").Value = 1499.5
$ws.Range("I78").Value = 1499
$ws.Range("J78").Value = 1500
$ws.Range("K78").Value = 13491
$ws.Range("L78").Value = 13500
$ws.Range("M78").Value = -8499
$ws.Range("N78").Value = -23484
$ws.Range("H126").Value = 3999
$ws.Range("I126").Value = 3999
$ws.Range("K126").Value = 11997
$ws.Range("M126").Value = -7057

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 44982.4
$ws.Range("J15").Value = 44982.4
$ws.Range("L15").Value = 44982.4
$ws.Range("N15").Value = -45558.4
$ws.Range("H57").Value = 58375
$ws.Range("I57").Value = 22666.666
$ws.Range("J57").Value = 79800
$ws.Range("K57").Value = 22666.666
$ws.Range("L57").Value = 79800
$ws.Range("M57").Value = -21846.666
$ws.Range("N57").Value = -81440
$ws.Range("H70").Value = 3999
$ws.Range("I70").Value = 3999
$ws.Range("K70").Value = 3999
$ws.Range("M70").Value = -3729
$ws.Range("H73").Value = 3999
$ws.Range("I73").Value = 3999
$ws.Range("K73").Value = 3999
$ws.Range("M73").Value = -3063
$ws.Range("H80").Value = 27699.75
$ws.Range("J80").Value = 100000
$ws.Range("L80").Value = 100000
$ws.Range("N80").Value = -101996
$ws.Range("H81").Value = 44982.4
$ws.Range("J81").Value = 44982.4
$ws.Range("L81").Value = 44982.4
$ws.Range("N81").Value = -46978.4
$ws.Range("H83").Value = 27699.75
$ws.Range("J83").Value = 100000
$ws.Range("L83").Value = 500000
$ws.Range("N83").Value = -509984
$ws.Range("H84").Value = 44982.4
$ws.Range("J84").Value = 44982.4
$ws.Range("L84").Value = 134947.2
$ws.Range("N84").Value = -144931.2
$ws.Range("H132").Value = 5658
$ws.Range("J132").Value = 2613.5
$ws.Range("L132").Value = 7840.5
$ws.Range("N132").Value = -12900.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4091.6155
$ws.Range("J40").Value = 4599
$ws.Range("L40").Value = 4599
$ws.Range("N40").Value = -4871
$ws.Range("H55").Value = 1598.6154
$ws.Range("I55").Value = 1511.3334
$ws.Range("J55").Value = 1673.4286
$ws.Range("K55").Value = 1511.3334
$ws.Range("L55").Value = 1673.4286
$ws.Range("M55").Value = -1338.3334
$ws.Range("N55").Value = -2019.4286
$ws.Range("H122").Value = 7958.227
$ws.Range("I122").Value = 8614.299999999999
$ws.Range("K122").Value = 25842.9
$ws.Range("M122").Value = -23392.9

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

